$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('K2').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNQ5MK545%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T180010Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQDawd0ZNM9TT1hun0zBicadGO0YHu4JPipxQbbuwvMbMAIgHVc5Wuebr%2FNCJ3aXPmFg2%2BHSgJgJXvWX%2BjeQKM5dV%2FQq5gMIkP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKNP8rOIoFi4BA%2B01Sq6A2uUy%2BqOcCWwnDCfszOZfWvhKQwfN09w0Gul8hmdiRIWEWaSig3J53%2B83TUCPC7idF2zCVY18Dumc9Y8f1w1ENdqXjV2nXqPJyY6%2BM9PkGBDHi4EFyyKsKmrWULgbCX2tiDoSrnFj1b43vb0DkFlpAe1YD3o9uvgePxRRUYyN7f4US0axiylJAlflD06rEZ4qmMnG6OBsHwkC9wK4zlz8sGzYYFSWxMbS55VAKNtVOEkfM5jZyjbhcy6iS4oguxRCUl9ZpJP%2FxOYkhZl%2FwWcBBNA5luW7rHeAlJk%2BfFsNRnXZs7R84PAdYYJd0ERL6fQ0Y0huEhlGhRB0G3zijMOARU92SPetuY9juGy14PhVVQPFiVlT9tE5gIMSEdC%2FrlVkwWr7PF0Odbf%2BiSst0j1IiNjsghjfi3JtGiqlOvFPiSUG0XBfTmvJzzTz1pfo%2BAsjAVRytqKpTCf%2Bv7U06R8XZ2Id8CEl9GWFnTOJAnmRVAyD7ZbEo4NHvsE%2BCIwhpcbnvS60KZng5bm%2BbjUwJ98MDBAjfh%2F7UjbC5AQ1TZWYaBVgpUSOULuQbJP0kP0BLtN6fNwJETYeRGFE%2BIwosalvwY6pQHckFNg1wQJ98GZpYEJ94k4MURXjd56Gj3ESXuuXD5QVN6zcGUs7f6ePEpppuhHSNWQnw%2BpacASpBvigrAn%2BBQ3yDgA%2BupFRaFwkOCse3eFZQuZJitESDT2UNoCfj8vCdLhipjjdoTahjctJ3khc%2F2soZ1SOLQfGGnmeDbYwWRH1M9SROKKZTioOo97BYsNElGeDTbk2G2GtmIBJSXlrz4149YiHYs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2002%3A00%3A10%20GMT&X-Amz-Signature=92000388b34166c8933b30cec4e99f7f1b35aaaf914dd639b860e5b76b4566bc'
$ws.Range('K3').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNQ5MK545%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T180010Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQDawd0ZNM9TT1hun0zBicadGO0YHu4JPipxQbbuwvMbMAIgHVc5Wuebr%2FNCJ3aXPmFg2%2BHSgJgJXvWX%2BjeQKM5dV%2FQq5gMIkP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKNP8rOIoFi4BA%2B01Sq6A2uUy%2BqOcCWwnDCfszOZfWvhKQwfN09w0Gul8hmdiRIWEWaSig3J53%2B83TUCPC7idF2zCVY18Dumc9Y8f1w1ENdqXjV2nXqPJyY6%2BM9PkGBDHi4EFyyKsKmrWULgbCX2tiDoSrnFj1b43vb0DkFlpAe1YD3o9uvgePxRRUYyN7f4US0axiylJAlflD06rEZ4qmMnG6OBsHwkC9wK4zlz8sGzYYFSWxMbS55VAKNtVOEkfM5jZyjbhcy6iS4oguxRCUl9ZpJP%2FxOYkhZl%2FwWcBBNA5luW7rHeAlJk%2BfFsNRnXZs7R84PAdYYJd0ERL6fQ0Y0huEhlGhRB0G3zijMOARU92SPetuY9juGy14PhVVQPFiVlT9tE5gIMSEdC%2FrlVkwWr7PF0Odbf%2BiSst0j1IiNjsghjfi3JtGiqlOvFPiSUG0XBfTmvJzzTz1pfo%2BAsjAVRytqKpTCf%2Bv7U06R8XZ2Id8CEl9GWFnTOJAnmRVAyD7ZbEo4NHvsE%2BCIwhpcbnvS60KZng5bm%2BbjUwJ98MDBAjfh%2F7UjbC5AQ1TZWYaBVgpUSOULuQbJP0kP0BLtN6fNwJETYeRGFE%2BIwosalvwY6pQHckFNg1wQJ98GZpYEJ94k4MURXjd56Gj3ESXuuXD5QVN6zcGUs7f6ePEpppuhHSNWQnw%2BpacASpBvigrAn%2BBQ3yDgA%2BupFRaFwkOCse3eFZQuZJitESDT2UNoCfj8vCdLhipjjdoTahjctJ3khc%2F2soZ1SOLQfGGnmeDbYwWRH1M9SROKKZTioOo97BYsNElGeDTbk2G2GtmIBJSXlrz4149YiHYs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2002%3A00%3A10%20GMT&X-Amz-Signature=4e8234ef2db03b26947fbbf5172a51d7d2232ee3da0fe6dfe0392c1dec1cf385'
$ws.Range('K4').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNQ5MK545%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T180010Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQDawd0ZNM9TT1hun0zBicadGO0YHu4JPipxQbbuwvMbMAIgHVc5Wuebr%2FNCJ3aXPmFg2%2BHSgJgJXvWX%2BjeQKM5dV%2FQq5gMIkP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKNP8rOIoFi4BA%2B01Sq6A2uUy%2BqOcCWwnDCfszOZfWvhKQwfN09w0Gul8hmdiRIWEWaSig3J53%2B83TUCPC7idF2zCVY18Dumc9Y8f1w1ENdqXjV2nXqPJyY6%2BM9PkGBDHi4EFyyKsKmrWULgbCX2tiDoSrnFj1b43vb0DkFlpAe1YD3o9uvgePxRRUYyN7f4US0axiylJAlflD06rEZ4qmMnG6OBsHwkC9wK4zlz8sGzYYFSWxMbS55VAKNtVOEkfM5jZyjbhcy6iS4oguxRCUl9ZpJP%2FxOYkhZl%2FwWcBBNA5luW7rHeAlJk%2BfFsNRnXZs7R84PAdYYJd0ERL6fQ0Y0huEhlGhRB0G3zijMOARU92SPetuY9juGy14PhVVQPFiVlT9tE5gIMSEdC%2FrlVkwWr7PF0Odbf%2BiSst0j1IiNjsghjfi3JtGiqlOvFPiSUG0XBfTmvJzzTz1pfo%2BAsjAVRytqKpTCf%2Bv7U06R8XZ2Id8CEl9GWFnTOJAnmRVAyD7ZbEo4NHvsE%2BCIwhpcbnvS60KZng5bm%2BbjUwJ98MDBAjfh%2F7UjbC5AQ1TZWYaBVgpUSOULuQbJP0kP0BLtN6fNwJETYeRGFE%2BIwosalvwY6pQHckFNg1wQJ98GZpYEJ94k4MURXjd56Gj3ESXuuXD5QVN6zcGUs7f6ePEpppuhHSNWQnw%2BpacASpBvigrAn%2BBQ3yDgA%2BupFRaFwkOCse3eFZQuZJitESDT2UNoCfj8vCdLhipjjdoTahjctJ3khc%2F2soZ1SOLQfGGnmeDbYwWRH1M9SROKKZTioOo97BYsNElGeDTbk2G2GtmIBJSXlrz4149YiHYs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2002%3A00%3A10%20GMT&X-Amz-Signature=6811cc56f1cb3d4e9be6ac76e6060f94cf11d3b041cd54533395a0dae24aec7f'
$ws.Range('K5').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNQ5MK545%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T180010Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQDawd0ZNM9TT1hun0zBicadGO0YHu4JPipxQbbuwvMbMAIgHVc5Wuebr%2FNCJ3aXPmFg2%2BHSgJgJXvWX%2BjeQKM5dV%2FQq5gMIkP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKNP8rOIoFi4BA%2B01Sq6A2uUy%2BqOcCWwnDCfszOZfWvhKQwfN09w0Gul8hmdiRIWEWaSig3J53%2B83TUCPC7idF2zCVY18Dumc9Y8f1w1ENdqXjV2nXqPJyY6%2BM9PkGBDHi4EFyyKsKmrWULgbCX2tiDoSrnFj1b43vb0DkFlpAe1YD3o9uvgePxRRUYyN7f4US0axiylJAlflD06rEZ4qmMnG6OBsHwkC9wK4zlz8sGzYYFSWxMbS55VAKNtVOEkfM5jZyjbhcy6iS4oguxRCUl9ZpJP%2FxOYkhZl%2FwWcBBNA5luW7rHeAlJk%2BfFsNRnXZs7R84PAdYYJd0ERL6fQ0Y0huEhlGhRB0G3zijMOARU92SPetuY9juGy14PhVVQPFiVlT9tE5gIMSEdC%2FrlVkwWr7PF0Odbf%2BiSst0j1IiNjsghjfi3JtGiqlOvFPiSUG0XBfTmvJzzTz1pfo%2BAsjAVRytqKpTCf%2Bv7U06R8XZ2Id8CEl9GWFnTOJAnmRVAyD7ZbEo4NHvsE%2BCIwhpcbnvS60KZng5bm%2BbjUwJ98MDBAjfh%2F7UjbC5AQ1TZWYaBVgpUSOULuQbJP0kP0BLtN6fNwJETYeRGFE%2BIwosalvwY6pQHckFNg1wQJ98GZpYEJ94k4MURXjd56Gj3ESXuuXD5QVN6zcGUs7f6ePEpppuhHSNWQnw%2BpacASpBvigrAn%2BBQ3yDgA%2BupFRaFwkOCse3eFZQuZJitESDT2UNoCfj8vCdLhipjjdoTahjctJ3khc%2F2soZ1SOLQfGGnmeDbYwWRH1M9SROKKZTioOo97BYsNElGeDTbk2G2GtmIBJSXlrz4149YiHYs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2002%3A00%3A10%20GMT&X-Amz-Signature=7425865fc777f0b95584707c0add5b7dd495f4d6fb1bee9d6dd49ef04efb7076'
$ws.Range('K6').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNQ5MK545%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T180010Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIQDawd0ZNM9TT1hun0zBicadGO0YHu4JPipxQbbuwvMbMAIgHVc5Wuebr%2FNCJ3aXPmFg2%2BHSgJgJXvWX%2BjeQKM5dV%2FQq5gMIkP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKNP8rOIoFi4BA%2B01Sq6A2uUy%2BqOcCWwnDCfszOZfWvhKQwfN09w0Gul8hmdiRIWEWaSig3J53%2B83TUCPC7idF2zCVY18Dumc9Y8f1w1ENdqXjV2nXqPJyY6%2BM9PkGBDHi4EFyyKsKmrWULgbCX2tiDoSrnFj1b43vb0DkFlpAe1YD3o9uvgePxRRUYyN7f4US0axiylJAlflD06rEZ4qmMnG6OBsHwkC9wK4zlz8sGzYYFSWxMbS55VAKNtVOEkfM5jZyjbhcy6iS4oguxRCUl9ZpJP%2FxOYkhZl%2FwWcBBNA5luW7rHeAlJk%2BfFsNRnXZs7R84PAdYYJd0ERL6fQ0Y0huEhlGhRB0G3zijMOARU92SPetuY9juGy14PhVVQPFiVlT9tE5gIMSEdC%2FrlVkwWr7PF0Odbf%2BiSst0j1IiNjsghjfi3JtGiqlOvFPiSUG0XBfTmvJzzTz1pfo%2BAsjAVRytqKpTCf%2Bv7U06R8XZ2Id8CEl9GWFnTOJAnmRVAyD7ZbEo4NHvsE%2BCIwhpcbnvS60KZng5bm%2BbjUwJ98MDBAjfh%2F7UjbC5AQ1TZWYaBVgpUSOULuQbJP0kP0BLtN6fNwJETYeRGFE%2BIwosalvwY6pQHckFNg1wQJ98GZpYEJ94k4MURXjd56Gj3ESXuuXD5QVN6zcGUs7f6ePEpppuhHSNWQnw%2BpacASpBvigrAn%2BBQ3yDgA%2BupFRaFwkOCse3eFZQuZJitESDT2UNoCfj8vCdLhipjjdoTahjctJ3khc%2F2soZ1SOLQfGGnmeDbYwWRH1M9SROKKZTioOo97BYsNElGeDTbk2G2GtmIBJSXlrz4149YiHYs%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2002%3A00%3A10%20GMT&X-Amz-Signature=45d6f3aa5ff073bb58ef15587b8cac2276d689938f3814dd2801ca5175c105b6'
